# Fix mojibake "Â±" (double-encoded UTF-8) -> "±" (plus-minus sign)
# in the f1_score_weighted / training_time / test_time columns (B, C, D)
# for data rows 2 through 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c1 = [char]194
$c2 = [char]177
$badChar = "$c1$c2"     # "Â±" (mojibake: U+00C2 U+00B1)
$goodChar = [char]177   # "±" (U+00B1 PLUS-MINUS SIGN)

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $value = [string]$cell.Value2
        if ($value.Contains($badChar)) {
            $cell.Value2 = $value.Replace($badChar, $goodChar)
        }
    }
}
